# "wrapping up test file audit"
#
# The optimization_parameters sheet had picked up a stray leftover row
# (row 16: label "Sheet" with junk values 3 / 4) that doesn't belong with
# the rest of the parameter table. Remove it -- this also shifts the
# "simulation_timepoints" row up from 17 to 16.
#
# The rest of the change is just the view/selection state left behind by
# the audit pass: the degradation_rates sheet was scrolled to D40, the
# optimization_parameters sheet was left on A42 after the row delete, and
# the workbook was finally parked on the threshold_b tab.

$wb = $excel.ActiveWorkbook

$wsDegradation = $wb.Worksheets.Item("degradation_rates")
$wsDegradation.Range("D40").Select()

$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Rows.Item(16).Delete()
$wsOpt.Range("A42").Select()

$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
